# project.xlsx update:
#  - reassign the owner of the "3- Presentation..." task from Denise to Vladi
#  - add two new tasks: "Create ppt skeleton" (Denise/presentation/done) and
#    "turn all Denise's labs into notebooks" (Denise/code/tbd)
#  - refresh the autofilter (owner=Denise, status=tbd) and the hidden rows that
#    result from it
#  - move the active selection to B17

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- data changes -----------------------------------------------------
# Row 7 ("3- Presentation - Prompt engineering...") owner: Denise -> Vladi
$ws.Range("C7").Value = "Vladi"

# New row 15
$ws.Range("A15").Value = "Create ppt skeleton"
$ws.Range("B15").Value = "presentation"
$ws.Range("C15").Value = "Denise"
$ws.Range("D15").Value = "done"

# New row 16
$ws.Range("A16").Value = "turn all Denise's labs into notebooks"
$ws.Range("B16").Value = "code"
$ws.Range("C16").Value = "Denise"
$ws.Range("D16").Value = "tbd"

# Re-fit the rows we just wrote into so the emulator doesn't stamp a
# stray custom row height on them.
$ws.Rows.Item(7).AutoFit()
$ws.Rows.Item(15).AutoFit()
$ws.Rows.Item(16).AutoFit()

# --- autofilter ---------------------------------------------------------
$ws.AutoFilterMode = $false
$ws.Range("A1:D16").AutoFilter() | Out-Null
$ws.Range("A1:D16").AutoFilter(3, @("Denise"), 7) | Out-Null
$ws.Range("A1:D16").AutoFilter(4, @("tbd"), 7) | Out-Null

# Keep the _FilterDatabase defined name range in sync with the filter.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "project!_FilterDatabase") {
        $n.RefersTo = "=project!`$A`$1:`$D`$16"
    }
}

# --- row visibility -------------------------------------------------
# Rows now hidden by the new filter (owner=Denise AND status=tbd)
$ws.Rows.Item(2).Hidden = $true
$ws.Rows.Item(3).Hidden = $true
$ws.Rows.Item(4).Hidden = $true
$ws.Rows.Item(5).Hidden = $true
$ws.Rows.Item(6).Hidden = $true
$ws.Rows.Item(7).Hidden = $true
$ws.Rows.Item(8).Hidden = $true
$ws.Rows.Item(9).Hidden = $true
$ws.Rows.Item(10).Hidden = $true
$ws.Rows.Item(11).Hidden = $true
$ws.Rows.Item(12).Hidden = $true
$ws.Rows.Item(13).Hidden = $true
$ws.Rows.Item(14).Hidden = $false
$ws.Rows.Item(15).Hidden = $true
$ws.Rows.Item(16).Hidden = $false

# --- selection ---------------------------------------------------------
$ws.Range("B17").Select() | Out-Null
